# Add "Bare Point" (garda) data to the inventory matrix.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the empty column B (spacer column) so the data columns
#    (spawning/hatching/water temps/experimental model/literature model)
#    close up from C:G to B:F.
$ws.Columns("B").Delete()

# 2. Insert a new row under the header for the "garda" lake entry,
#    pushing the existing lake rows down by one.
$ws.Rows("2").Insert()

# 3. Rewrite the full grid explicitly so every cell ends up exactly right.
$headers = @("lake", "spawning", "hatching", "water temps", "experimental model", "literature model")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

$data = @(
    @("garda",            "?? Checking", "",                   "?? Checking", "",   "Eckmann, 1987"),
    @("geneva",            "x",           "Dan?",               "x",            "x",  "Eckmann, 1987"),
    @("bourget",            "",            "",                   "x",            "x",  "Eckmann, 1987"),
    @("annecy",             "",            "x",                  "x",            "",   "Eckmann, 1987"),
    @("constance",          "",            "",                   "",             "x",  "Eckmann, 1987"),
    @("konnevesi",          "x",           "x",                  "x",            "x",  "Luczyński & Kirklewska, 1984"),
    @("ontario",            "x",           "x",                  "x",            "x",  "Colby & Brooke, 1973"),
    @("apostle islands",    "x",           "x",                  "x",            "x",  "Colby & Brooke, 1973"),
    @("thunder bay",        "x",           "Jared? Dan (GP)?",   "in process",   "x",  "Colby & Brooke, 1973")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    $excelRow = $r + 2
    for ($c = 0; $c -lt $row.Length; $c++) {
        $value = $row[$c]
        if ($value -ne "") {
            $ws.Cells.Item($excelRow, $c + 1).Value = $value
        } else {
            $ws.Cells.Item($excelRow, $c + 1).Value = $null
        }
    }
}

# 4. Column widths: A, C, D, E already carry their correct (bestFit) widths
#    forward from the deleted column B's old neighbours - only the new
#    literature-model column (F) needs an explicit width. (ColumnWidth adds
#    a fixed ~5/6-character pixel pad when stored, so back it out here to
#    land on a rendered width of exactly 25.)
$ws.Columns("F").ColumnWidth = 145/6

# 5. Selection reflects where the editor last clicked when finishing this edit.
$ws.Range("H8").Select()
